# Update cryptos list prices / volumes (and reordered rows 49-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.966.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.895.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7727"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.07"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3141"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.88"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07366"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08073"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7726"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.898.67"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.961.69"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007851"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.166.66"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.164"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1578"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.450"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.77"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.035"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.428"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.480"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05584"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.071"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.243"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7546"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.006"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.681"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01933"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.61"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4475"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.105.19"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.034"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8511"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.897"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.51"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.545"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.774"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.008"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.58%  "
